$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1333333333333333
$ws.Range("C2").Value = 0.6566666666666666
$ws.Range("P2").Value = 0.09333333333333334
$ws.Range("S2").Value = 0.1166666666666667
$ws.Range("B3").Value = 0.005025125628140704
$ws.Range("C3").Value = 0.01507537688442211
$ws.Range("J3").Value = 0.03517587939698492
$ws.Range("P3").Value = 0.7537688442211056
$ws.Range("S3").Value = 0.1909547738693467
$ws.Range("J4").Value = 0.03225806451612903
$ws.Range("O4").Value = 0.01612903225806452
$ws.Range("P4").Value = 0.7096774193548387
$ws.Range("S4").Value = 0.2419354838709677
$ws.Range("B6").Value = 0.04680851063829787
$ws.Range("D6").Value = 0.01276595744680851
$ws.Range("E6").Value = 0.00425531914893617
$ws.Range("F6").Value = 0.06808510638297872
$ws.Range("J6").Value = 0.2723404255319149
$ws.Range("O6").Value = 0.01276595744680851
$ws.Range("Q6").Value = 0.1574468085106383
$ws.Range("R6").Value = 0.0851063829787234
$ws.Range("S6").Value = 0.3404255319148936
$ws.Range("B7").Value = 0.1228070175438596
$ws.Range("D7").Value = 0.04093567251461988
$ws.Range("F7").Value = 0.02339181286549707
$ws.Range("J7").Value = 0.1286549707602339
$ws.Range("O7").Value = 0.02923976608187134
$ws.Range("Q7").Value = 0.1695906432748538
$ws.Range("R7").Value = 0.09941520467836257
$ws.Range("S7").Value = 0.3859649122807017
$ws.Range("B8").Value = 0.1108870967741935
$ws.Range("D8").Value = 0.02016129032258064
$ws.Range("F8").Value = 0.07862903225806452
$ws.Range("J8").Value = 0.1068548387096774
$ws.Range("O8").Value = 0.03024193548387097
$ws.Range("Q8").Value = 0.1834677419354839
$ws.Range("R8").Value = 0.1129032258064516
$ws.Range("S8").Value = 0.3568548387096774
$ws.Range("B9").Value = 0.1116751269035533
$ws.Range("D9").Value = 0.03045685279187817
$ws.Range("E9").Value = 0.005076142131979695
$ws.Range("F9").Value = 0.04568527918781726
$ws.Range("J9").Value = 0.1065989847715736
$ws.Range("O9").Value = 0.03045685279187817
$ws.Range("Q9").Value = 0.1624365482233502
$ws.Range("R9").Value = 0.1015228426395939
$ws.Range("S9").Value = 0.4060913705583756
$ws.Range("B10").Value = 0.1067251461988304
$ws.Range("D10").Value = 0.02704678362573099
$ws.Range("E10").Value = 0.0007309941520467836
$ws.Range("F10").Value = 0.07383040935672515
$ws.Range("J10").Value = 0.1067251461988304
$ws.Range("O10").Value = 0.02046783625730994
$ws.Range("Q10").Value = 0.2010233918128655
$ws.Range("R10").Value = 0.1067251461988304
$ws.Range("S10").Value = 0.3567251461988304
$ws.Range("G11").Value = 0.1338289962825279
$ws.Range("J11").Value = 0.1152416356877323
$ws.Range("K11").Value = 0.1895910780669145
$ws.Range("L11").Value = 0.5539033457249071
$ws.Range("S11").Value = 0.007434944237918215
$ws.Range("G12").Value = 0.6687898089171974
$ws.Range("J12").Value = 0.2292993630573248
$ws.Range("L12").Value = 0.03821656050955414
$ws.Range("S12").Value = 0.06369426751592357
$ws.Range("G13").Value = 0.7083333333333334
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.04166666666666666
$ws.Range("F15").Value = 0.0299625468164794
$ws.Range("H15").Value = 0.1947565543071161
$ws.Range("I15").Value = 0.06741573033707865
$ws.Range("J15").Value = 0.3258426966292135
$ws.Range("K15").Value = 0.06367041198501873
$ws.Range("M15").Value = 0.003745318352059925
$ws.Range("O15").Value = 0.0449438202247191
$ws.Range("S15").Value = 0.2696629213483146
$ws.Range("F16").Value = 0.01382488479262673
$ws.Range("H16").Value = 0.1889400921658986
$ws.Range("I16").Value = 0.1013824884792627
$ws.Range("J16").Value = 0.4009216589861751
$ws.Range("K16").Value = 0.07373271889400922
$ws.Range("M16").Value = 0.01382488479262673
$ws.Range("O16").Value = 0.05529953917050692
$ws.Range("S16").Value = 0.152073732718894
$ws.Range("F17").Value = 0.01948051948051948
$ws.Range("H17").Value = 0.1861471861471861
$ws.Range("I17").Value = 0.09956709956709957
$ws.Range("J17").Value = 0.4329004329004329
$ws.Range("K17").Value = 0.0735930735930736
$ws.Range("M17").Value = 0.02597402597402598
$ws.Range("O17").Value = 0.04761904761904762
$ws.Range("S17").Value = 0.1147186147186147
$ws.Range("F18").Value = 0.03891050583657588
$ws.Range("H18").Value = 0.1750972762645914
$ws.Range("I18").Value = 0.07392996108949416
$ws.Range("J18").Value = 0.4513618677042802
$ws.Range("K18").Value = 0.07392996108949416
$ws.Range("M18").Value = 0.007782101167315175
$ws.Range("O18").Value = 0.07392996108949416
$ws.Range("S18").Value = 0.1050583657587549
$ws.Range("F19").Value = 0.01381427475057559
$ws.Range("H19").Value = 0.2110514198004605
$ws.Range("I19").Value = 0.07214121258633922
$ws.Range("J19").Value = 0.3844973138910207
$ws.Range("K19").Value = 0.1005372217958557
$ws.Range("M19").Value = 0.02455871066768995
$ws.Range("N19").Value = 0.001534919416730622
$ws.Range("O19").Value = 0.08135072908672294
$ws.Range("S19").Value = 0.1105141980046048

Write-Host "Applied 109 cell updates"
